# feat: add 2022-Q1 data
#
# Before: sheets are "2021-Q2", "2021-Q4", "总计" (the per-date roll-up).
# After:  the old "总计" sheet becomes "2022-Q1" (holding the new quarter's
#         per-fund breakdown) and a fresh "总计" sheet is appended at the
#         end with the roll-up table, now including the "2022-Q1" row.

$wb = $excel.ActiveWorkbook

$wsQ4  = $wb.Worksheets.Item(2)   # "2021-Q4"  -- style/layout template
$wsQ1  = $wb.Worksheets.Item(3)   # currently "总计", about to become "2022-Q1"

# ---------------------------------------------------------------------------
# Step 1: rename the existing "总计" sheet to "2022-Q1".
# ---------------------------------------------------------------------------
$wsQ1.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# Step 2: insert a brand-new "总计" sheet right after it (sheetId 4, last
# tab), matching the workbook.xml ordering in the diff.
# ---------------------------------------------------------------------------
$wsNewTot = $wb.Worksheets.Add($null, $wsQ1)
$wsNewTot.Name = "总计"

# A neutral, unstyled helper cell used to strip the transient "quote
# prefix" formatting that Excel stamps on a cell when a value is entered
# with a leading apostrophe (our trick for forcing numeric-looking text to
# stay text). Copying its (plain) format back onto a cell removes that
# flag again without touching the cell's value.
$plainFmt = $wsQ1.Cells.Item(500, 500)
$plainFmt.Value = 0
$plainFmt.Copy()

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.PasteSpecial(-4122)   # xlPasteFormats
}

# ===========================================================================
# "2022-Q1" sheet: per-fund holdings table, columns A-H.
# Seed the block from the "2021-Q4" sheet via Copy so the header/index
# styling (cellXf index 2) transfers correctly, then overwrite the text.
# ===========================================================================

$wsQ4.Range("B1:H3").Copy($wsQ1.Range("B1"))

$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

Set-TextValue $wsQ1.Cells.Item(2,2) "007178"
$wsQ1.Cells.Item(2,3).Value = "浙商港股通中华交易服务预期高股息指数增强A"
Set-TextValue $wsQ1.Cells.Item(2,4) "7.93"
Set-TextValue $wsQ1.Cells.Item(2,5) "90.20"
Set-TextValue $wsQ1.Cells.Item(2,6) "7.38"
Set-TextValue $wsQ1.Cells.Item(2,7) "0.5852"
$wsQ1.Cells.Item(2,8).Value = 6

Set-TextValue $wsQ1.Cells.Item(3,2) "007216"
$wsQ1.Cells.Item(3,3).Value = "浙商港股通中华交易服务预期高股息指数增强C"
Set-TextValue $wsQ1.Cells.Item(3,4) "4.60"
Set-TextValue $wsQ1.Cells.Item(3,5) "90.20"
Set-TextValue $wsQ1.Cells.Item(3,6) "7.38"
Set-TextValue $wsQ1.Cells.Item(3,7) "0.3395"
$wsQ1.Cells.Item(3,8).Value = 6

# ===========================================================================
# "总计" sheet: date roll-up table, columns A-D, with the new "2022-Q1"
# row inserted on top and the previous rows shifted down.
# ===========================================================================

$wsQ4.Range("B1:D1").Copy($wsNewTot.Range("B1"))
$wsNewTot.Cells.Item(1,2).Value = "日期"
$wsNewTot.Cells.Item(1,3).Value = "持有数量(只)"
$wsNewTot.Cells.Item(1,4).Value = "持有市值(亿元)"

$wsQ4.Range("A2").Copy($wsNewTot.Range("A2"))
$wsQ4.Range("A2").Copy($wsNewTot.Range("A3"))
$wsQ4.Range("A2").Copy($wsNewTot.Range("A4"))

$wsNewTot.Cells.Item(2,1).Value = 0
$wsNewTot.Cells.Item(2,2).Value = "2022-Q1"
$wsNewTot.Cells.Item(2,3).Value = 2
$wsNewTot.Cells.Item(2,4).Value = 0.92

$wsNewTot.Cells.Item(3,1).Value = 1
$wsNewTot.Cells.Item(3,2).Value = "2021-Q4"
$wsNewTot.Cells.Item(3,3).Value = 2
$wsNewTot.Cells.Item(3,4).Value = 0.23

$wsNewTot.Cells.Item(4,1).Value = 2
$wsNewTot.Cells.Item(4,2).Value = "2021-Q2"
$wsNewTot.Cells.Item(4,3).Value = 1
$wsNewTot.Cells.Item(4,4).Value = 0.01

# Tidy up the helper cell used for format-stripping.
$plainFmt.ClearContents()

Write-Host "2022-Q1 data added"
